$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DSD_LIVESTOCK")

# Row 4: REF_AREA concept ID renamed to GEO_PICT
$ws.Range("A4").Value = "GEO_PICT"

# Rows 8-11: codelist for the four new livestock dimensions consolidated
# onto the shared CL_COM_YESNO codelist instead of per-dataset codelists.
$ws.Range("F8").Value = "CL_COM_YESNO"
$ws.Range("F9").Value = "CL_COM_YESNO"
$ws.Range("F10").Value = "CL_COM_YESNO"
$ws.Range("F11").Value = "CL_COM_YESNO"

# Update the active selection to reflect the final edit position.
$ws.Activate()
$ws.Range("F7").Select()
